$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so they stay as text (matching source data)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2").Value = "56.219.87"
$ws.Range("E2").Value = "  -4.99%  "
$ws.Range("D3").Value = "2.353.78"
$ws.Range("E3").Value = "  -6.62%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.12"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.94"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  -2.79%  "
$ws.Range("D9").Value = "2.366.78"
$ws.Range("E9").Value = "  -6.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0954"
$ws.Range("E10").Value = "  -4.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  -8.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.314"
$ws.Range("E13").Value = "  -5.79%  "
$ws.Range("D14").Value = "2.774.88"
$ws.Range("E14").Value = "  -6.50%  "
$ws.Range("D15").Value = "56.163.72"
$ws.Range("E15").Value = "  -5.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.28"
$ws.Range("E16").Value = "  -5.20%  "
$ws.Range("E17").Value = "  -4.95%  "
$ws.Range("D18").Value = "2.391.51"
$ws.Range("E18").Value = "  -5.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.20"
$ws.Range("E19").Value = "  -4.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.05"
$ws.Range("E20").Value = "  -4.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.06"
$ws.Range("E21").Value = "  -3.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.25"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.385"
$ws.Range("E26").Value = "  -6.51%  "
$ws.Range("D27").Value = "2.468.77"
$ws.Range("E27").Value = "  -6.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.152"
$ws.Range("E28").Value = "  -5.25%  "
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.47"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("D32").Value = "0.0₃0710"
$ws.Range("E32").Value = "  -7.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.11"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("E34").Value = "  -7.90%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.62"
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("E38").Value = "  -6.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.68"
$ws.Range("E39").Value = "  -7.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.57"
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.793"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -7.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("E44").Value = "  -6.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "122.80"
$ws.Range("E45").Value = "  -7.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.564"
$ws.Range("E46").Value = "  -5.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "251.43"
$ws.Range("E47").Value = "  -9.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0900"
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0484"
$ws.Range("E49").Value = "  -5.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0205"
$ws.Range("E50").Value = "  -6.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.47"
$ws.Range("E51").Value = "  -7.37%  "
